# "Contenu du stage" statistics (rows 16-23) are refreshed with the real
# soutenance counts and their corresponding percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the percentage cells so values such as "3.85 %"
# are not auto-converted to numeric percentages by Excel's input parser.
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G20").NumberFormat = "@"

# Row 16: C# -> 1 soutenance, 3.85 %
$ws.Range("E16").Value = 1
$ws.Range("G16").Value = "3.85 %"

# Row 17: COBOL -> 23 soutenances, 88.46 %
$ws.Range("E17").Value = 23
$ws.Range("G17").Value = "88.46 %"

# Row 20: ANDROID -> 2 soutenances, 7.69 %
$ws.Range("E20").Value = 2
$ws.Range("G20").Value = "7.69 %"

# Drop the temporary "@" text format again so the cells keep the workbook's
# default (unstyled) look, matching the original formatting.
$ws.Range("G16").ClearFormats()
$ws.Range("G17").ClearFormats()
$ws.Range("G20").ClearFormats()
